$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    "even_MAG-GUT1026.fa",
    "even_MAG-GUT10417.fa",
    "even_MAG-GUT1139.fa",
    "even_MAG-GUT11412.fa",
    "even_MAG-GUT1148.fa",
    "even_MAG-GUT12217.fa",
    "even_MAG-GUT13955.fa",
    "even_MAG-GUT18024.fa",
    "even_MAG-GUT19599.fa",
    "even_MAG-GUT22049.fa",
    "even_MAG-GUT22878.fa",
    "even_MAG-GUT28136.fa",
    "even_MAG-GUT29051.fa",
    "even_MAG-GUT29076.fa",
    "even_MAG-GUT39136.fa",
    "even_MAG-GUT40857.fa",
    "even_MAG-GUT42485.fa",
    "even_MAG-GUT42494.fa",
    "even_MAG-GUT42584.fa",
    "even_MAG-GUT4338.fa",
    "even_MAG-GUT43894.fa",
    "even_MAG-GUT456.fa",
    "even_MAG-GUT49243.fa",
    "even_MAG-GUT49507.fa",
    "even_MAG-GUT52107.fa",
    "even_MAG-GUT52138.fa",
    "even_MAG-GUT5727.fa",
    "even_MAG-GUT59039.fa",
    "even_MAG-GUT61159.fa",
    "even_MAG-GUT61959.fa",
    "even_MAG-GUT6280.fa",
    "even_MAG-GUT6290.fa",
    "even_MAG-GUT70200.fa",
    "even_MAG-GUT75471.fa",
    "even_MAG-GUT77982.fa",
    "even_MAG-GUT78879.fa",
    "even_MAG-GUT78908.fa",
    "even_MAG-GUT83946.fa",
    "even_MAG-GUT86868.fa",
    "even_MAG-GUT87091.fa",
    "even_MAG-GUT87486.fa",
    "even_MAG-GUT87573.fa",
    "even_MAG-GUT87828.fa",
    "even_MAG-GUT88085.fa",
    "even_MAG-GUT88218.fa",
    "even_MAG-GUT88257.fa",
    "even_MAG-GUT88679.fa",
    "even_MAG-GUT88862.fa"
)

$values = @(
    70146.62000423382,
    95809.76437652647,
    71509.47446510397,
    114550.1212106485,
    79840.37822910654,
    108368.8770634211,
    115384.0047471105,
    104884.1258428711,
    117405.1599745854,
    97421.89571959394,
    84300.51216472185,
    81207.03595207828,
    73972.36771533091,
    83485.85692499307,
    112241.7689173326,
    66602.25207459994,
    77062.31740954959,
    102461.0152130505,
    81549.14162279168,
    111757.8563361279,
    81832.4187264478,
    101048.0196982883,
    73076.030053158,
    71610.60477488833,
    116802.7749086313,
    110298.9152929693,
    84473.23288990806,
    77457.33608977369,
    87885.79209551282,
    89277.35330399006,
    62958.46883102979,
    62958.46883102979,
    114961.8156116052,
    110012.6219574911,
    105234.9571069992,
    113890.5585103906,
    93748.30563322347,
    91274.61671756732,
    107439.9103055472,
    117380.4233117023,
    118380.9714551817,
    118498.23652985,
    118957.9554701441,
    123162.7889049486,
    120740.43272986,
    112605.7526683244,
    112176.3614426853,
    113033.5730537835
)

# Extend column A formatting (bordered/bold/centered style, matching existing data rows)
# to the newly added rows before writing data into them.
$ws.Range("A2").Copy()
$ws.Range("A32:A49").PasteSpecial(-4122)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
    $ws.Cells.Item($row, 3).Value = "o__Christensenellales"
    $ws.Cells.Item($row, 4).Value = "o__Christensenellales"
}
